# GlycoNet Tech Workshops Budget template edit
#
# The "Honoraria (non-GlycoNet only):" line item (row 20) in the
# FACILITATORS section is removed entirely. Deleting the whole row lets
# Excel naturally:
#   - shift every row below it up by one (old row 21 -> 20, 22 -> 21, 23 -> 22)
#   - re-point the FACILITATOR TOTAL SUM() formulas to the new (smaller)
#     Travel:/Accommodations: range
#   - re-point the WORKSHOP TOTAL SUM() formula and its merged cell range
#   - drop the now-unused "Honoraria (non-GlycoNet only):" shared string
#   - shrink the sheet dimension by one row
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Honoraria (non-GlycoNet only):" row.
$ws.Rows.Item(20).Delete() | Out-Null

# Match the author's resulting selection.
$ws.Range("A20").Select() | Out-Null
